# CV.xlsx — "Add files via upload"
# The author re-saved the workbook after the referenced images on disk were
# renamed to lowercase (FHS.jpg -> fhs.jpg, HSLU.jpg -> hslu.jpg,
# Migros.jpg -> migros.jpg). The only semantic change is that the four
# cells in column H ("Bild") that referenced those files now reference the
# lowercase filenames instead. (Excel's shared-string table reorders
# itself automatically as a result: the old casing becomes unused and is
# dropped, while the new casing is appended at the end — matching the
# target diff.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 3: Hochschule für Angewandte Wissenschaften (FHS) St.Gallen
$ws.Range("H4").Value = "hslu.jpg"
$ws.Range("H3").Value = "fhs.jpg"

# Rows 8 & 9: Genossenschaft Migros Ostschweiz, Gossau SG
$ws.Range("H8").Value = "migros.jpg"
$ws.Range("H9").Value = "migros.jpg"

# The workbook was also left scrolled/selected differently on save.
$ws.Range("D6").Select()
